$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to be treated as text so numeric-looking values
# (e.g. "1.00", "26.00", "64.812.83") keep their exact literal formatting
# instead of being auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Cells.Item(2, 4).Value = '64.812.83'
$ws.Cells.Item(2, 5).Value = '  +1.77%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '3.160.13'
$ws.Cells.Item(3, 5).Value = '  +0.91%  '

# Row 4
$ws.Cells.Item(4, 4).Value = '1.00'
$ws.Cells.Item(4, 5).Value = '  +0.41%  '

# Row 5
$ws.Cells.Item(5, 4).Value = '592.97'
$ws.Cells.Item(5, 5).Value = '  +1.12%  '

# Row 6
$ws.Cells.Item(6, 4).Value = '153.09'
$ws.Cells.Item(6, 5).Value = '  +4.85%  '

# Row 7
$ws.Cells.Item(7, 5).Value = '  +0.18%  '

# Row 8
$ws.Cells.Item(8, 4).Value = '3.156.86'
$ws.Cells.Item(8, 5).Value = '  +0.84%  '

# Row 9
$ws.Cells.Item(9, 4).Value = '0.536'
$ws.Cells.Item(9, 5).Value = '  +1.49%  '

# Row 10
$ws.Cells.Item(10, 4).Value = '0.162'
$ws.Cells.Item(10, 5).Value = '  +1.06%  '

# Row 11
$ws.Cells.Item(11, 4).Value = '6.00'
$ws.Cells.Item(11, 5).Value = '  +4.55%  '

# Row 12
$ws.Cells.Item(12, 4).Value = '0.466'
$ws.Cells.Item(12, 5).Value = '  +1.94%  '

# Row 13
$ws.Cells.Item(13, 4).Value = '38.70'
$ws.Cells.Item(13, 5).Value = '  +5.30%  '

# Row 14
$ws.Cells.Item(14, 4).Value = '0.0000249'
$ws.Cells.Item(14, 5).Value = '  +1.33%  '

# Row 15
$ws.Cells.Item(15, 2).Value = 'WrappedliquidstakedEther2.0'
$ws.Cells.Item(15, 3).Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Cells.Item(15, 4).Value = '3.676.63'
$ws.Cells.Item(15, 5).Value = '  +0.62%  '

# Row 16
$ws.Cells.Item(16, 2).Value = 'TRON'
$ws.Cells.Item(16, 3).Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Cells.Item(16, 4).Value = '0.121'
$ws.Cells.Item(16, 5).Value = '  -0.17%  '

# Row 17
$ws.Cells.Item(17, 2).Value = 'Polkadot'
$ws.Cells.Item(17, 3).Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Cells.Item(17, 4).Value = '7.32'
$ws.Cells.Item(17, 5).Value = '  +3.81%  '

# Row 18
$ws.Cells.Item(18, 2).Value = 'WrappedBTC'
$ws.Cells.Item(18, 3).Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Cells.Item(18, 4).Value = '64.416.56'
$ws.Cells.Item(18, 5).Value = '  +1.32%  '

# Row 19
$ws.Cells.Item(19, 2).Value = 'WrappedEther'
$ws.Cells.Item(19, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Cells.Item(19, 4).Value = '3.156.41'
$ws.Cells.Item(19, 5).Value = '  +0.64%  '

# Row 20
$ws.Cells.Item(20, 4).Value = '475.78'
$ws.Cells.Item(20, 5).Value = '  +2.67%  '

# Row 21
$ws.Cells.Item(21, 4).Value = '15.01'
$ws.Cells.Item(21, 5).Value = '  +5.39%  '

# Row 22
$ws.Cells.Item(22, 4).Value = '0.759'
$ws.Cells.Item(22, 5).Value = '  +3.72%  '

# Row 23
$ws.Cells.Item(23, 4).Value = '7.71'
$ws.Cells.Item(23, 5).Value = '  +4.17%  '

# Row 24
$ws.Cells.Item(24, 4).Value = '13.52'
$ws.Cells.Item(24, 5).Value = '  +4.64%  '

# Row 25
$ws.Cells.Item(25, 4).Value = '2.43'
$ws.Cells.Item(25, 5).Value = '  +10.75%  '

# Row 26
$ws.Cells.Item(26, 4).Value = '82.46'
$ws.Cells.Item(26, 5).Value = '  +1.79%  '

# Row 27
$ws.Cells.Item(27, 5).Value = '  +0.07%  '

# Row 28
$ws.Cells.Item(28, 4).Value = '10.02'
$ws.Cells.Item(28, 5).Value = '  +8.63%  '

# Row 29
$ws.Cells.Item(29, 4).Value = '2.74'
$ws.Cells.Item(29, 5).Value = '  +2.37%  '

# Row 30
$ws.Cells.Item(30, 4).Value = '7.42'
$ws.Cells.Item(30, 5).Value = '  +6.63%  '

# Row 31
$ws.Cells.Item(31, 2).Value = 'ImmutableX'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(31, 4).Value = '2.24'
$ws.Cells.Item(31, 5).Value = '  +1.27%  '

# Row 32
$ws.Cells.Item(32, 2).Value = 'FirstDigitalUSD'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Cells.Item(32, 4).Value = '1.01'
$ws.Cells.Item(32, 5).Value = '  +0.68%  '

# Row 33
$ws.Cells.Item(33, 5).Value = '  +7.19%  '

# Row 34
$ws.Cells.Item(34, 4).Value = '27.87'
$ws.Cells.Item(34, 5).Value = '  +3.66%  '

# Row 35
$ws.Cells.Item(35, 4).Value = '0.0₃0881'
$ws.Cells.Item(35, 5).Value = '  +4.94%  '

# Row 36
$ws.Cells.Item(36, 2).Value = 'dogwifhat'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Cells.Item(36, 4).Value = '3.55'
$ws.Cells.Item(36, 5).Value = '  +7.77%  '

# Row 37
$ws.Cells.Item(37, 2).Value = 'Mantle'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Cells.Item(37, 4).Value = '1.07'
$ws.Cells.Item(37, 5).Value = '  +3.03%  '

# Row 38
$ws.Cells.Item(38, 4).Value = '6.23'
$ws.Cells.Item(38, 5).Value = '  +3.82%  '

# Row 39
$ws.Cells.Item(39, 4).Value = '2.33'
$ws.Cells.Item(39, 5).Value = '  +1.90%  '

# Row 40
$ws.Cells.Item(40, 4).Value = '468.58'
$ws.Cells.Item(40, 5).Value = '  +6.98%  '

# Row 41
$ws.Cells.Item(41, 4).Value = '9.38'
$ws.Cells.Item(41, 5).Value = '  +6.87%  '

# Row 42
$ws.Cells.Item(42, 4).Value = '51.43'
$ws.Cells.Item(42, 5).Value = '  +0.68%  '

# Row 43
$ws.Cells.Item(43, 4).Value = '0.301'
$ws.Cells.Item(43, 5).Value = '  +8.69%  '

# Row 44
$ws.Cells.Item(44, 4).Value = '0.0377'
$ws.Cells.Item(44, 5).Value = '  +1.89%  '

# Row 45
$ws.Cells.Item(45, 4).Value = '2.899.48'
$ws.Cells.Item(45, 5).Value = '  -0.29%  '

# Row 46
$ws.Cells.Item(46, 4).Value = '0.111'
$ws.Cells.Item(46, 5).Value = '  +3.53%  '

# Row 47
$ws.Cells.Item(47, 4).Value = '38.33'
$ws.Cells.Item(47, 5).Value = '  +3.71%  '

# Row 48
$ws.Cells.Item(48, 4).Value = '131.41'
$ws.Cells.Item(48, 5).Value = '  +3.88%  '

# Row 49
$ws.Cells.Item(49, 4).Value = '26.00'
$ws.Cells.Item(49, 5).Value = '  +7.79%  '

# Row 50
$ws.Cells.Item(50, 4).Value = '2.31'
$ws.Cells.Item(50, 5).Value = '  +6.87%  '

# Row 51
$ws.Cells.Item(51, 5).Value = '  +0.03%  '
